# Actualizar ligas y agregar Liga Argentina actualizada
# Append the new Liga Suecia 2025 fixtures (rows 124-131) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Fecha, Local, Visita, GolesLocal, GolesVisita, FixtureID,
#             CornersLocal, CornersVisita, AmarillasLocal, AmarillasVisita,
#             RojasLocal, RojasVisita, Goles1TLocal, Goles1TVisita,
#             Goles2TLocal, Goles2TVisita, PosesionLocal, PosesionVisita, Resultado
$nuevosPartidos = @(
    @("2025-07-19", "Djurgardens IF", "IF Elfsborg",       1, 0, 1342066, 7,  1,  3, 2, 0, 0, 0, 0, 1, 0, "46%", "54%", "L"),
    @("2025-07-19", "Osters IF",      "Malmo FF",          0, 2, 1342068, 1,  5,  0, 1, 0, 0, 0, 0, 0, 2, "37%", "63%", "V"),
    @("2025-07-19", "Degerfors IF",   "Gais",              0, 3, 1342073, 6,  1,  1, 1, 0, 0, 0, 0, 0, 3, "53%", "47%", "V"),
    @("2025-07-20", "Sirius",         "IFK Goteborg",      0, 1, 1342071, 3,  6,  5, 1, 0, 0, 0, 0, 0, 1, "59%", "41%", "V"),
    @("2025-07-20", "Mjallby AIF",    "AIK Stockholm",     2, 0, 1342069, 5,  5,  2, 3, 0, 0, 0, 0, 2, 0, "55%", "45%", "L"),
    @("2025-07-20", "Hammarby FF",    "IF Brommapojkarna", 3, 2, 1342067, 8,  6,  1, 1, 0, 0, 0, 0, 3, 2, "70%", "30%", "L"),
    @("2025-07-20", "Halmstad",       "BK Hacken",         0, 0, 1342072, 2, 11,  3, 2, 0, 0, 0, 0, 0, 0, "33%", "67%", "E"),
    @("2025-07-21", "IFK Norrkoping", "IFK Varnamo",       3, 1, 1342070, 4,  7,  2, 1, 0, 0, 0, 0, 3, 1, "51%", "49%", "L")
)

$startRow = 124

for ($i = 0; $i -lt $nuevosPartidos.Count; $i++) {
    $fila = $startRow + $i
    $datos = $nuevosPartidos[$i]

    # Column A (Fecha) holds an ISO-looking date string that must stay text,
    # not get auto-converted into a date serial number.
    $cFecha = $ws.Cells.Item($fila, 1)
    $cFecha.NumberFormat = "@"
    $cFecha.Value = $datos[0]
    $cFecha.Style = "Normal"

    $ws.Cells.Item($fila, 2).Value = $datos[1]
    $ws.Cells.Item($fila, 3).Value = $datos[2]
    $ws.Cells.Item($fila, 4).Value = $datos[3]
    $ws.Cells.Item($fila, 5).Value = $datos[4]
    $ws.Cells.Item($fila, 6).Value = $datos[5]
    $ws.Cells.Item($fila, 7).Value = $datos[6]
    $ws.Cells.Item($fila, 8).Value = $datos[7]
    $ws.Cells.Item($fila, 9).Value = $datos[8]
    $ws.Cells.Item($fila, 10).Value = $datos[9]
    $ws.Cells.Item($fila, 11).Value = $datos[10]
    $ws.Cells.Item($fila, 12).Value = $datos[11]
    $ws.Cells.Item($fila, 13).Value = $datos[12]
    $ws.Cells.Item($fila, 14).Value = $datos[13]
    $ws.Cells.Item($fila, 15).Value = $datos[14]
    $ws.Cells.Item($fila, 16).Value = $datos[15]

    # Columns Q & R (possession %) are plain text like "46%", not a numeric
    # percentage value, so they need the same text-forcing treatment.
    $cPosLocal = $ws.Cells.Item($fila, 17)
    $cPosLocal.NumberFormat = "@"
    $cPosLocal.Value = $datos[16]
    $cPosLocal.Style = "Normal"

    $cPosVisita = $ws.Cells.Item($fila, 18)
    $cPosVisita.NumberFormat = "@"
    $cPosVisita.Value = $datos[17]
    $cPosVisita.Style = "Normal"

    $ws.Cells.Item($fila, 19).Value = $datos[18]
}
